$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Learning goals")

# --- Learning goal 1 planning rows (columns C/D/F), 9.00 - 12.00 slot (rows 14-16) ---
$ws.Range("C14").Value = 44579
$ws.Range("C14").NumberFormat = "d-mmm"
$ws.Range("D14").Value = "9.00 - 12.00"
$ws.Range("F14").Value = 3

$ws.Range("C15").Value = 44580
$ws.Range("C15").NumberFormat = "d-mmm"
$ws.Range("D15").Value = "9.00 - 12.00"
$ws.Range("F15").Value = 3

$ws.Range("C16").Value = 44581
$ws.Range("C16").NumberFormat = "d-mmm"
$ws.Range("D16").Value = "9.00 - 12.00"
$ws.Range("F16").Value = 3

# --- Learning goal 2 planning rows (columns Q/R/T), rows 19-24 ---
$ws.Range("Q19").Value = 44580
$ws.Range("Q19").NumberFormat = "d-mmm"
$ws.Range("R19").Value = "9.00  - 12.00"
$ws.Range("T19").Value = 3

$ws.Range("Q20").Value = 44581
$ws.Range("Q20").NumberFormat = "d-mmm"
$ws.Range("R20").Value = "9.00 - 12.00"
$ws.Range("T20").Value = 3

$ws.Range("Q21").Value = 44588
$ws.Range("Q21").NumberFormat = "d-mmm"
$ws.Range("R21").Value = "9.00 - 12.00"
$ws.Range("T21").Value = 3

$ws.Range("Q22").Value = 44589
$ws.Range("Q22").NumberFormat = "d-mmm"
$ws.Range("R22").Value = "9.00 - 12.00"
$ws.Range("T22").Value = 3

$ws.Range("Q23").Value = 44592
$ws.Range("Q23").NumberFormat = "d-mmm"
$ws.Range("R23").Value = "9.00 - 12.00"
$ws.Range("T23").Value = 3

$ws.Range("Q24").Value = 44593
$ws.Range("Q24").NumberFormat = "d-mmm"
$ws.Range("R24").Value = "9.00 - 12.00"
$ws.Range("T24").Value = 3

# --- Learning goal 1 planning rows (columns C/D/F), 13.00 - 17.00 slot (rows 17-21) ---
$ws.Range("C17").Value = 44587
$ws.Range("C17").NumberFormat = "d-mmm"
$ws.Range("D17").Value = "13.00 - 17.00"
$ws.Range("F17").Value = 4

$ws.Range("C18").Value = 44588
$ws.Range("C18").NumberFormat = "d-mmm"
$ws.Range("D18").Value = "13.00 - 17.00"
$ws.Range("F18").Value = 4

$ws.Range("C19").Value = 44589
$ws.Range("C19").NumberFormat = "d-mmm"
$ws.Range("D19").Value = "13.00 - 17.00"
$ws.Range("F19").Value = 4

$ws.Range("C20").Value = 44592
$ws.Range("C20").NumberFormat = "d-mmm"
$ws.Range("D20").Value = "13.00 - 17.00"
$ws.Range("F20").Value = 4

$ws.Range("C21").Value = 44593
$ws.Range("C21").NumberFormat = "d-mmm"
$ws.Range("D21").Value = "13.00 - 17.00"
$ws.Range("F21").Value = 4

# --- Move the running-total formulas down to reflect the new rows ---
$ws.Range("F22").ClearContents()
$ws.Range("F23").Formula = "=SUM(F3:F21)"
$ws.Range("T26").Formula = "=SUM(T3:T24)"

# --- Final selection state ---
$ws.Range("S32").Select()
